{"js": "// Replace every occurrence of \"SAIL\" with \"SUN\" in the document (the intro\n// paragraph mentions it three times), then relocate the \"_GoBack\" bookmark\n// from the last list item (\"When you have completed ...\") to sit right\n// before \" Manual\" in the (now updated) intro paragraph.\n\n// Step 1: \"SAIL\" -> \"SUN\"\nconst sailHits = context.document.body.search(\"SAIL\", { matchCase: true });\nsailHits.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < sailHits.items.length; i++) {\n  sailHits.items[i].insertText(\"SUN\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Step 2: move the \"_GoBack\" bookmark.\n// Remove it from wherever it currently lives (\"When you have completed ...\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert it immediately before the word \"Manual\" in the intro paragraph.\nconst manualHits = context.document.body.search(\"Manual\", { matchCase: true });\nmanualHits.load(\"text\");\nawait context.sync();\n\nconst insertionPoint = manualHits.items[0].getRange(Word.RangeLocation.before);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Update the Quick Start Guide intro paragraph: rebrand \"SAIL\" -> \"SUN\", and\n# relocate the \"_GoBack\" bookmark from the last list item (\"When you have\n# completed ...\") to sit right before \" Manual\" in the intro paragraph.\n\n$d = $word.ActiveDocument\n\n# Step 1: replace every \"SAIL\" with \"SUN\" (all three occurrences live in the\n# intro paragraph: \"...with SAIL...\", \"...on SAIL...\", \"...SAIL Manual...\").\n$d.Content.Find.Execute(\"SAIL\", $false, $false, $false, $false, $false, $true, 1, $false, \"SUN\", 2) | Out-Null\n\n# Step 2: remove the \"_GoBack\" bookmark from wherever it currently sits.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete() | Out-Null\n}\n\n# Step 3: re-create it immediately before the word \"Manual\" in the intro\n# paragraph (now reading \"...please refer to the SUN Manual - English.\").\n$findRange = $d.Content\n$findRange.Find.Execute(\"Manual\") | Out-Null\n$findRange.Collapse(1)  # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $findRange) | Out-Null\n"}
